# Update the Greeting cell for rule R10 from "Good Morning" to "GIT UPDATE".
# Setting the cell value causes the now-unused "Good Morning" shared string
# to be dropped and a new "GIT UPDATE" shared string to be appended, which
# is exactly the shared-strings-table churn captured by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"

# Match the saved selection/active cell for the sheet (E8).
$ws.Range("E8").Select()
